$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $found = $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
    if (-not $found) {
        throw "Find/Replace failed for: $old"
    }
}

# 1. Update activation date
Replace-Text 'Ativação: 01/01/2018' 'Ativação: 01/01/2025'

# 2. Fix double space before PMBOX
Replace-Text 'Projects: Life cycle of Project. Feasibility Economic. Flowchart.  PMBOX. Complementary aspects of the preparation of a project.' 'Projects: Life cycle of Project. Feasibility Economic. Flowchart. PMBOX. Complementary aspects of the preparation of a project.'

# 3. Update Portuguese program summary (add spaces, new/changed items)
Replace-Text '1.Introdução: Conceitos de Gestão de Projetos2.Aspectos da Implantação de Projetos: Etapas Fundamentais e Formas Parciais3.O Ciclo de Vida do Projeto4.Aspectos da Viabilidade de Projetos: Receitas, Custos, Ponto de Nivelamento, Estimativas5.Guia PMBOK: Principais Áreas de Conhecimento6.Plano de Projeto7.O Gerente de Projeto e as Interfaces: Equipes de Projeto8.Legalização de Projetos: Aspectos sobre o Licenciamento9.Técnicas de Análise de Riscos Operacionais10.Aspectos sobre Auditorias e Auditorias Integradas11.Apresentação de Trabalhos.' '1. Introdução: Conceitos de Gestão de Projetos2. Aspectos da Implantação de Projetos: Etapas Fundamentais e Formas Parciais3. O Ciclo de Vida do Projeto4. Aspectos da Viabilidade de Projetos: Receitas, Custos, Ponto de Nivelamento, Estimativas5. Guia PMBOK: Principais Áreas de Conhecimento6. Plano de Projeto7. O Gerente de Projeto e as Interfaces: Equipes de Projeto8. Legalização de Projetos: Aspectos sobre o Licenciamento9. Gestão de Riscos - Técnicas de Análise de Riscos10. Aspectos sobre Auditorias e Auditorias Integradas11. Provas e/ou apresentações de Trabalhos12. Sistema de Gestão Integrada (SGI): Qualidade; Meio Ambiente; e Saúde e Segurança.13. Legislação brasileira aplicada ao SGI'

# 4. Update English program summary (add spaces)
Replace-Text '1. Introduction: Concepts of Project Management2. Aspects of Implementation of Projects: Fundamental steps and Partial Forms3. The Life cycle of Project4. Aspects of Feasibility of Project: Revenue, Costs, Point of Leveling, Estimates5. PMBOK Guide: Main areas of knowledge6. Plan of Project 7. The Manager of Project and Interfaces: Team of Project8. Legalization of Project: Aspects about the Licensing 9. Techniques of risk analysis operational10. Aspects about Audits and Integrated Audits11. Work Presentation.' '1. Introduction: Concepts of Project Management 2. Aspects of Implementation of Projects: Fundamental steps and Partial Forms 3. The Life cycle of Project 4. Aspects of Feasibility of Project: Revenue, Costs, Point of Leveling, Estimates 5. PMBOK Guide: Main areas of knowledge 6. Plan of Project 7. The Manager of Project and Interfaces: Team of Project 8. Legalization of Project: Aspects about the Licensing 9. Techniques of risk analysis operational 10. Aspects about Audits and Integrated Audits 11. Work Presentation.'

# 5. Merge bibliography list items (remove line breaks, update content)
$nl = [char]11
$biblioOld = '1) Chemical Engineering Plant Design (Vilbrandt e Bryden)' + $nl + '2) Project Engineering of Process Plants Rose e Barrow (2ª impressão - 1968)' + $nl + '3) Elaboração e Análise de Projetos Simonsen, M.H. e H. Flanger' + $nl + '4) Implantação de Indústrias. Valle, E.C. Livros Técnicos e Científicos Editores S/A, Rio de Janeiro. ' + $nl + '5) Introdução ao Projeto de Engenharia. Asimov, Morris. Editora Mestre Jou-São Paulo, 1968.' + $nl + '6) IBGR  Instituto Brasileiro de Gerenciamento de Riscos, 2000.'
$biblioNew = '1) Chemical Engineering Plant Design (Vilbrandt e Bryden) 2) Project Engineering of Process Plants Rose e Barrow (2ª impressão - 1968) 3) Elaboração e Análise de Projetos Simonsen, M.H. e H. Flanger 4) Implantação de Indústrias. Valle, E.C. Livros Técnicos e Científicos Editores S/A, Rio de Janeiro. 5) Gestão Integrada: Qualidade, Meio Ambiente, Prevenção. Antecipação de riscos e outras ferramentas para implantação. CHAVES. F. J. M., 1ª Ed. 2022.6) PMBOK Guide (6ª Ed., 2017)7) Normas ABNT NBR ISO: 9001; 14001; 45001. 8) IBGR – Instituto Brasileiro de Gerenciamento de Riscos, 2000.'
Replace-Text $biblioOld $biblioNew

